# Change UNDER exception logic to show current allocation period.
# Rewrite the Exceptions sheet rows 2-14 (EmployeeName, ExceptionType,
# StartDate, EndDate, FreeOrExcessPercent, SourceProjectsOrClients)
# so that StartDate/EndDate reflect the busy (allocated) period instead
# of the free period, per the new AllocationAnalyzer behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chen Noam",            "UNDER", "2026-01-18", "2026-02-15", 100, "None"),
    @("Galandor Moshe Yakov", "UNDER", "2026-01-18", "2026-04-18", 90,  "None"),
    @("Meir Zipora",          "UNDER", "2026-01-18", "2026-02-28", 87,  "None"),
    @("Noiman Yehuda",        "UNDER", "2026-01-18", "2026-01-31", 50,  "None"),
    @("Noiman Yehuda",        "UNDER", "2026-01-18", "2026-01-31", 100, "None"),
    @("Weingarten Ayala",     "UNDER", "2026-01-18", "2026-02-28", 85,  "None"),
    @("Itzhaki Yair",         "UNDER", "2026-01-18", "2026-04-18", 50,  "None"),
    @("Peretz Yehonathan",    "UNDER", "2026-01-18", "2026-03-31", 100, "None"),
    @("Levin Yanir",          "UNDER", "2026-01-18", "2026-02-28", 100, "None"),
    @("Halevy Maor",          "UNDER", "2026-01-18", "2026-02-28", 100, "None"),
    @("Cohen Aharon",         "UNDER", "2026-01-18", "2026-03-31", 100, "None"),
    @("Pruzanski Yossi",      "UNDER", "2026-01-18", "2026-03-31", 100, "None"),
    @("Morgenstern Elisheva", "UNDER", "2026-01-18", "2026-03-31", 100, "None")
)

# Keep the StartDate/EndDate columns formatted as Text so the date-like
# strings (e.g. "2026-01-18") are not auto-converted into Excel date
# serial numbers - they must remain plain text, matching the source data.
$ws.Range("C2:D14").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}
